# Refresh cached market-price derived figures (columns H-N) across the
# Leve profit sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 5124987.5
$ws.Range("I12").Value = 6223031
$ws.Range("J12").Value = 785
$ws.Range("K12").Value = 6223031
$ws.Range("L12").Value = 785
$ws.Range("M12").Value = -6222861
$ws.Range("N12").Value = -1125

$ws.Range("H53").Value = 631.1429000000001
$ws.Range("I53").Value = 156
$ws.Range("K53").Value = 156
$ws.Range("M53").Value = 481

$ws.Range("H100").Value = 2020.4117
$ws.Range("I100").Value = 2029
$ws.Range("K100").Value = 2029
$ws.Range("M100").Value = -1488

$ws.Range("H125").Value = 1897346.9
$ws.Range("I125").Value = 3790277
$ws.Range("J125").Value = 4416.6665
$ws.Range("K125").Value = 34112493
$ws.Range("L125").Value = 39749.9985
$ws.Range("M125").Value = -34110033
$ws.Range("N125").Value = -44669.9985

$ws.Range("H132").Value = 28818.621
$ws.Range("I132").Value = 32852.812
$ws.Range("J132").Value = 2999.8
$ws.Range("K132").Value = 98558.43599999999
$ws.Range("L132").Value = 8999.400000000001
$ws.Range("M132").Value = -96028.43599999999
$ws.Range("N132").Value = -14059.4

$ws.Range("H137").Value = 3031442.8
$ws.Range("I137").Value = 928.3333
$ws.Range("J137").Value = 8334843
$ws.Range("K137").Value = 2784.9999
$ws.Range("L137").Value = 25004529
$ws.Range("M137").Value = -234.9998999999998
$ws.Range("N137").Value = -25009629

$ws.Range("H138").Value = 7212.457
$ws.Range("J138").Value = 4394.7095
$ws.Range("L138").Value = 13184.1285
$ws.Range("N138").Value = -23464.1285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 187467.94
$ws.Range("I32").Value = 201375.14
$ws.Range("K32").Value = 201375.14
$ws.Range("M32").Value = -201088.14

$ws.Range("H43").Value = 17744.143
$ws.Range("J43").Value = 17492.2
$ws.Range("L43").Value = 17492.2
$ws.Range("N43").Value = -18118.2

$ws.Range("H61").Value = 1786700.4
$ws.Range("I61").Value = 3905.84
$ws.Range("J61").Value = 7357933.5
$ws.Range("K61").Value = 3905.84
$ws.Range("L61").Value = 7357933.5
$ws.Range("M61").Value = -3693.84
$ws.Range("N61").Value = -7358357.5

$ws.Range("H74").Value = 497038.5
$ws.Range("I74").Value = 1700.9375
$ws.Range("J74").Value = 1716331
$ws.Range("K74").Value = 1700.9375
$ws.Range("L74").Value = 1716331
$ws.Range("M74").Value = -826.9375
$ws.Range("N74").Value = -1718079

$ws.Range("H77").Value = 497038.5
$ws.Range("I77").Value = 1700.9375
$ws.Range("J77").Value = 1716331
$ws.Range("K77").Value = 8504.6875
$ws.Range("L77").Value = 8581655
$ws.Range("M77").Value = -4136.6875
$ws.Range("N77").Value = -8590391

$ws.Range("H136").Value = 1786700.4
$ws.Range("I136").Value = 3905.84
$ws.Range("J136").Value = 7357933.5
$ws.Range("K136").Value = 11717.52
$ws.Range("L136").Value = 22073800.5
$ws.Range("M136").Value = -9167.52
$ws.Range("N136").Value = -22078900.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1733.625
$ws.Range("J20").Value = 1962.25
$ws.Range("L20").Value = 1962.25
$ws.Range("N20").Value = -2456.25

$ws.Range("H105").Value = 9100.5625
$ws.Range("J105").Value = 4285.2856
$ws.Range("L105").Value = 4285.2856
$ws.Range("N105").Value = -7779.2856

$ws.Range("H107").Value = 6256.1924
$ws.Range("I107").Value = 6708.8223
$ws.Range("K107").Value = 6708.8223
$ws.Range("M107").Value = -4788.8223

$ws.Range("H134").Value = 64289216
$ws.Range("I134").Value = 3517.125
$ws.Range("J134").Value = 150003490
$ws.Range("K134").Value = 10551.375
$ws.Range("L134").Value = 450010470
$ws.Range("M134").Value = -8016.375
$ws.Range("N134").Value = -450015540

$ws.Range("H140").Value = 96663.336
$ws.Range("J140").Value = 96663.336
$ws.Range("L140").Value = 96663.336
$ws.Range("N140").Value = -107023.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3484.2
$ws.Range("I16").Value = 1874.5
$ws.Range("J16").Value = 5898.75
$ws.Range("K16").Value = 1874.5
$ws.Range("L16").Value = 5898.75
$ws.Range("M16").Value = -1587.5
$ws.Range("N16").Value = -6472.75

$ws.Range("H31").Value = 3409.2068
$ws.Range("I31").Value = 3314.0833
$ws.Range("K31").Value = 3314.0833
$ws.Range("M31").Value = -3019.0833

$ws.Range("H34").Value = 3409.2068
$ws.Range("I34").Value = 3314.0833
$ws.Range("K34").Value = 3314.0833
$ws.Range("M34").Value = -3112.0833

$ws.Range("H94").Value = 2139
$ws.Range("I94").Value = 2100
$ws.Range("K94").Value = 2100
$ws.Range("M94").Value = -1649

$ws.Range("H107").Value = 1356.3636
$ws.Range("I107").Value = 1141.375
$ws.Range("J107").Value = 1929.6666
$ws.Range("K107").Value = 1141.375
$ws.Range("L107").Value = 1929.6666
$ws.Range("M107").Value = 778.625
$ws.Range("N107").Value = -5769.6666

$ws.Range("H113").Value = 3484.2
$ws.Range("I113").Value = 1874.5
$ws.Range("J113").Value = 5898.75
$ws.Range("K113").Value = 1874.5
$ws.Range("L113").Value = 5898.75
$ws.Range("M113").Value = 295.5
$ws.Range("N113").Value = -10238.75

$ws.Range("H134").Value = 1993.7142
$ws.Range("I134").Value = 1718.2941
$ws.Range("J134").Value = 3164.25
$ws.Range("K134").Value = 5154.8823
$ws.Range("L134").Value = 9492.75
$ws.Range("M134").Value = -2619.8823
$ws.Range("N134").Value = -14562.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2711.0625
$ws.Range("I137").Value = 3025.4443
$ws.Range("J137").Value = 2306.8572
$ws.Range("K137").Value = 9076.332900000001
$ws.Range("L137").Value = 6920.571599999999
$ws.Range("M137").Value = -3976.332900000001
$ws.Range("N137").Value = -17120.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 13611
$ws.Range("J63").Value = 13611
$ws.Range("L63").Value = 13611
$ws.Range("N63").Value = -14983

$ws.Range("H66").Value = 13611
$ws.Range("J66").Value = 13611
$ws.Range("L66").Value = 40833
$ws.Range("N66").Value = -47697

$ws.Range("H70").Value = 4723.0625
$ws.Range("I70").Value = 4713.0835
$ws.Range("K70").Value = 4713.0835
$ws.Range("M70").Value = -4443.0835

$ws.Range("H73").Value = 4723.0625
$ws.Range("I73").Value = 4713.0835
$ws.Range("K73").Value = 4713.0835
$ws.Range("M73").Value = -3777.0835

$ws.Range("H126").Value = 3437.9
$ws.Range("I126").Value = 3351.8333
$ws.Range("J126").Value = 3567
$ws.Range("K126").Value = 10055.4999
$ws.Range("L126").Value = 10701
$ws.Range("M126").Value = -7585.499899999999
$ws.Range("N126").Value = -15641

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3463.6365
$ws.Range("I16").Value = 3211.7646
$ws.Range("K16").Value = 3211.7646
$ws.Range("M16").Value = -3041.7646

$ws.Range("H22").Value = 8676.23
$ws.Range("J22").Value = 8874.25
$ws.Range("L22").Value = 8874.25
$ws.Range("N22").Value = -9464.25

$ws.Range("H27").Value = 8676.23
$ws.Range("J27").Value = 8874.25
$ws.Range("L27").Value = 8874.25
$ws.Range("N27").Value = -9088.25

$ws.Range("H122").Value = 3568.3713
$ws.Range("I122").Value = 2937.96
$ws.Range("K122").Value = 8813.880000000001
$ws.Range("M122").Value = -6363.880000000001

$ws.Range("H132").Value = 3090.0571
$ws.Range("I132").Value = 2975.5806
$ws.Range("K132").Value = 8926.7418
$ws.Range("M132").Value = -6396.7418

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4538018.5
$ws.Range("I81").Value = 7561881.5
$ws.Range("J81").Value = 2224
$ws.Range("K81").Value = 15123763
$ws.Range("L81").Value = 4448
$ws.Range("M81").Value = -15122702
$ws.Range("N81").Value = -6570

$ws.Range("H84").Value = 4538018.5
$ws.Range("I84").Value = 7561881.5
$ws.Range("J84").Value = 2224
$ws.Range("K84").Value = 75618815
$ws.Range("L84").Value = 22240
$ws.Range("M84").Value = -75613511
$ws.Range("N84").Value = -32848

$ws.Range("H107").Value = 1242919.1
$ws.Range("I107").Value = 705.44446
$ws.Range("K107").Value = 2116.33338
$ws.Range("M107").Value = -196.33338

$ws.Range("H122").Value = 2997.9546
$ws.Range("I122").Value = 2817.2778
$ws.Range("K122").Value = 8451.8334
$ws.Range("M122").Value = -6001.8334

$ws.Range("H132").Value = 2188.8708
$ws.Range("I132").Value = 2160.5173
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 6481.5519
$ws.Range("L132").Value = 7800
$ws.Range("M132").Value = -3951.5519
$ws.Range("N132").Value = -12860
